# Update metrics_40_6 sheet: rotate model ids (model_40_6_0 -> 12,
# model_40_6_12 -> 24, model_40_6_24 -> 0) and refresh every row's metric
# values per the commit "atualizado todo o treinamento para o novo lm".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model id labels (column A) for the three rows whose id changed.
$ws.Cells.Item(2, 1).Value = "model_40_6_12"
$ws.Cells.Item(14, 1).Value = "model_40_6_24"
$ws.Cells.Item(26, 1).Value = "model_40_6_0"

# Refreshed metric values (columns B..Q) shared by rows 2-25.
$commonValues = @(
    0.9999845286516351,
    0.9991389381953254,
    0.9999458973752269,
    0.9999826861946171,
    0.9999540210781167,
    0.00001444182458493536,
    0.000803763398418822,
    0.000082863676070671,
    0.000007397727364258409,
    0.0000451308600559713,
    0.0001971186733916371,
    0.003800240069381849,
    1.00001125188972,
    0.003962024156716867,
    136.2907641523688,
    205.7666861698563
)

for ($row = 2; $row -le 25; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $commonValues[$col - 2]
    }
}

# Row 26 (model_40_6_0) carries its own freshly computed values.
$row26Values = @(
    0.9999845285973327,
    0.9991389379755921,
    0.9999458973705828,
    0.9999826861946171,
    0.9999540212394307,
    0.00001444187527390777,
    0.0008037636035302677,
    0.00008286368318347278,
    0.000007397727364258409,
    0.00004513070171746471,
    0.000197127007959379,
    0.003800246738556297,
    1.000011251929213,
    0.003962031109811879,
    136.2907571326348,
    205.7666791501223
)

for ($col = 2; $col -le 17; $col++) {
    $ws.Cells.Item(26, $col).Value = $row26Values[$col - 2]
}
